$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("nivel")
$ws.Range("B123").Value = 495004.884416534
$ws.Range("B125").Value = 497029.796576029
$ws.Range("B126").Value = 496425.468224716
$ws.Range("B127").Value = 491587.062858516
$ws.Range("B128").Value = 470868.305527377
$ws.Range("B129").Value = 445852.447872867
$ws.Range("B130").Value = 424257.644076143
$ws.Range("B131").Value = 419200.025924307
$ws.Range("B132").Value = 425268.865857971
$ws.Range("B133").Value = 433921.98098691
$ws.Range("B134").Value = 440611.641640823
$ws.Range("B135").Value = 448785.874236871
$ws.Range("B136").Value = 471279.481076438
$ws.Range("B137").Value = 489389.038898445
$ws.Range("B138").Value = 486756.28867032
$ws.Range("B139").Value = 478425.424719264
$ws.Range("B140").Value = 484354.845988484
$ws.Range("B141").Value = 494926.283216486
$ws.Range("B142").Value = 518512.745506591
$ws.Range("B143").Value = 533429.024796907
$ws.Range("B144").Value = 527363.202630966
$ws.Range("B145").Value = 535686.312722413
$ws.Range("B146").Value = 557628.523954013
$ws.Range("B147").Value = 564423.797294381
$ws.Range("B148").Value = 572994.4839706169
$ws.Range("B149").Value = 587500.9624074301
$ws.Range("B150").Value = 604350.836630969
$ws.Range("B151").Value = 621320.584365549
$ws.Range("B152").Value = 620705.4021628191
$ws.Range("B153").Value = 640163.5519535129
$ws.Range("B154").Value = 647511.572007653
$ws.Range("B155").Value = 656783.657528617
$ws.Range("B156").Value = 655185.263083306
$ws.Range("B157").Value = 628102.3416506069
$ws.Range("B158").Value = 611861.4572936289
$ws.Range("B159").Value = 600012.417984769
$ws.Range("B160").Value = 611185.484875494
$ws.Range("B161").Value = 629321.3142874771
$ws.Range("B162").Value = 646923.40274078
$ws.Range("B163").Value = 674259.720502361
$ws.Range("B164").Value = 674453.4931644599
$ws.Range("B165").Value = 681549.969081769
$ws.Range("B166").Value = 698616.679696331
$ws.Range("B167").Value = 711708.808490338
$ws.Range("B168").Value = 716692.843817329
$ws.Range("B169").Value = 714948.488109337
$ws.Range("B170").Value = 709886.663286911
$ws.Range("B171").Value = 701246.370049731
$ws.Range("B172").Value = 708119.275846552
$ws.Range("B173").Value = 718422.846887617
$ws.Range("B174").Value = 715256.203557977
$ws.Range("B175").Value = 722085.858874657
$ws.Range("B176").Value = 725965.845637177
$ws.Range("B177").Value = 718947.844352816
$ws.Range("B178").Value = 708732.2925000919
$ws.Range("B179").Value = 704955.298531755
$ws.Range("B180").Value = 696621.877304511
$ws.Range("B181").Value = 700196.164114941
$ws.Range("B182").Value = 713047.970930157
$ws.Range("B183").Value = 723337.204341816
$ws.Range("B184").Value = 725727.972634828
$ws.Range("B185").Value = 717642.257969926
$ws.Range("B186").Value = 710078.1063833541
$ws.Range("B187").Value = 707251.193400725
$ws.Range("B188").Value = 702624.605382216
$ws.Range("B189").Value = 706411.390000167
$ws.Range("B190").Value = 710961.500410893
$ws.Range("B191").Value = 728468.988674912
$ws.Range("B192").Value = 729969.575873518
$ws.Range("B193").Value = 735797.957875601
$ws.Range("B194").Value = 734026.84932093
$ws.Range("B195").Value = 702748.636350124
$ws.Range("B196").Value = 704313.852217745
$ws.Range("B197").Value = 687618.714776894
$ws.Range("B198").Value = 691513.297529099
$ws.Range("B199").Value = 709920.5806747261
$ws.Range("B200").Value = 692652.7494384981
$ws.Range("B201").Value = 678094.223353247
$ws.Range("B202").Value = 636454.7497215091
$ws.Range("B203").Value = 611620.7201199461
$ws.Range("B204").Value = 621663.086153583
$ws.Range("B205").Value = 649643.15760734
$ws.Range("B206").Value = 677236.63131827
$ws.Range("B207").Value = 684764.787040495
$ws.Range("B208").Value = 692986.8406457281
$ws.Range("B209").Value = 704321.6656675281
$ws.Range("B210").Value = 722121.266143609
$ws.Range("B211").Value = 744773.66071154
$ws.Range("B212").Value = 736036.0666005909
$ws.Range("B213").Value = 722141.192043945
$ws.Range("B214").Value = 723656.571741228
$ws.Range("B215").Value = 718843.746122845
$ws.Range("B216").Value = 714697.4114862659
$ws.Range("B217").Value = 709098.621870624
$ws.Range("B218").Value = 695437.223225022
$ws.Range("B219").Value = 695275.3287403811
$ws.Range("B220").Value = 714021.332474257
$ws.Range("B221").Value = 726943.904959862
$ws.Range("B222").Value = 735155.8717923349
$ws.Range("B223").Value = 742728.137340707
$ws.Range("B224").Value = 736786.087496091

$ws = $wb.Worksheets.Item("trimestrales")
$ws.Range("B123").Value = -1.951321405003736
$ws.Range("B124").Value = -0.001650019048926588
$ws.Range("B125").Value = 0.4107259231463045
$ws.Range("B126").Value = -0.1215879521662688
$ws.Range("B127").Value = -0.9746489001667857
$ws.Range("B128").Value = -4.214666922001986
$ws.Range("B129").Value = -5.312707897485691
$ws.Range("B130").Value = -4.843486651189521
$ws.Range("B131").Value = -1.192110082742148
$ws.Range("B132").Value = 1.447719360294086
$ws.Range("B133").Value = 2.034739860742341
$ws.Range("B134").Value = 1.54167360655435
$ws.Range("B135").Value = 1.855201230182546
$ws.Range("B136").Value = 5.012102236465399
$ws.Range("B137").Value = 3.842636598700078
$ws.Range("B138").Value = -0.5379667337975125
$ws.Range("B139").Value = -1.711506177724709
$ws.Range("B140").Value = 1.239361656563154
$ws.Range("B141").Value = 2.182581079874923
$ws.Range("B142").Value = 4.765651590943709
$ws.Range("B143").Value = 2.876743034685236
$ws.Range("B144").Value = -1.137137629181395
$ws.Range("B145").Value = 1.578250065594977
$ws.Range("B146").Value = 4.096093312537219
$ws.Range("B147").Value = 1.218602178415162
$ws.Range("B148").Value = 1.518484287395472
$ws.Range("B149").Value = 2.531696001031136
$ws.Range("B150").Value = 2.868058999340595
$ws.Range("B151").Value = 2.807929882116178
$ws.Range("B152").Value = -0.09901204276986375
$ws.Range("B153").Value = 3.134844601463582
$ws.Range("B154").Value = 1.147834804358494
$ws.Range("B155").Value = 1.431956728157813
$ws.Range("B156").Value = -0.2433669636856428
$ws.Range("B157").Value = -4.13362799176017
$ws.Range("B158").Value = -2.585706704149227
$ws.Range("B159").Value = -1.936555925792471
$ws.Range("B160").Value = 1.862139275092245
$ws.Range("B161").Value = 2.967320046168553
$ws.Range("B162").Value = 2.796995438368755
$ws.Range("B163").Value = 4.225588013320714
$ws.Range("B164").Value = 0.02873857894916121
$ws.Range("B165").Value = 1.052181653624951
$ws.Range("B166").Value = 2.504102617384829
$ws.Range("B167").Value = 1.874007474270112
$ws.Range("B168").Value = 0.7002913646050102
$ws.Range("B169").Value = -0.2433895807722908
$ws.Range("B170").Value = -0.7079985350849416
$ws.Range("B171").Value = -1.217137000035162
$ws.Range("B172").Value = 0.9800985916452731
$ws.Range("B173").Value = 1.455061511882616
$ws.Range("B174").Value = -0.4407770915636489
$ws.Range("B175").Value = 0.9548543980054225
$ws.Range("B176").Value = 0.5373303901238025
$ws.Range("B177").Value = -0.9667123221480667
$ws.Range("B178").Value = -1.420903050612787
$ws.Range("B179").Value = -0.5329225164855123
$ws.Range("B180").Value = -1.182120518081131
$ws.Range("B181").Value = 0.5130885099762006
$ws.Range("B182").Value = 1.835458043598526
$ws.Range("B183").Value = 1.442993154897687
$ws.Range("B184").Value = 0.3305191933528961
$ws.Range("B185").Value = -1.11415226776308
$ws.Range("B186").Value = -1.054028173866106
$ws.Range("B187").Value = -0.3981129621116408
$ws.Range("B188").Value = -0.6541647524499261
$ws.Range("B189").Value = 0.5389484781693543
$ws.Range("B190").Value = 0.6441162295988789
$ws.Range("B191").Value = 2.462508624433357
$ws.Range("B192").Value = 0.2059919120696696
$ws.Range("B193").Value = 0.7984417700023227
$ws.Range("B194").Value = -0.2407058263364292
$ws.Range("B195").Value = -4.261181045317663
$ws.Range("B196").Value = 0.2227276990177218
$ws.Range("B197").Value = -2.370411626618074
$ws.Range("B198").Value = 0.5663869624998208
$ws.Range("B199").Value = 2.66188419100537
$ws.Range("B200").Value = -2.432360986043847
$ws.Range("B201").Value = -2.101850616640577
$ws.Range("B202").Value = -6.14066190179684
$ws.Range("B203").Value = -3.901931694661642
$ws.Range("B204").Value = 1.64192704780628
$ws.Range("B205").Value = 4.50084170621714
$ws.Range("B206").Value = 4.247481619379756
$ws.Range("B207").Value = 1.111599014892484
$ws.Range("B208").Value = 1.200712092800238
$ws.Range("B209").Value = 1.635647945527818
$ws.Range("B210").Value = 2.5271976347925
$ws.Range("B211").Value = 3.136923897685917
$ws.Range("B212").Value = -1.173187851810087
$ws.Range("B213").Value = -1.887798055986567
$ws.Range("B214").Value = 0.2098453479705054
$ws.Range("B215").Value = -0.6650703947595815
$ws.Range("B216").Value = -0.5768061082735532
$ws.Range("B217").Value = -0.7833790252575357
$ws.Range("B218").Value = -1.926586545826714
$ws.Range("B219").Value = -0.02327952534524336
$ws.Range("B220").Value = 2.696198607800127
$ws.Range("B221").Value = 1.809830028582637
$ws.Range("B222").Value = 1.129656191687367
$ws.Range("B223").Value = 1.030021773465628
$ws.Range("B224").Value = -0.8000302595093656

$ws = $wb.Worksheets.Item("i.a.")
$ws.Range("B123").Value = -0.8043457378072327
$ws.Range("B125").Value = -1.819654266967807
$ws.Range("B126").Value = -1.669937585145331
$ws.Range("B127").Value = -0.6904621884785245
$ws.Range("B128").Value = -4.874458839464236
$ws.Range("B129").Value = -10.29663594732465
$ws.Range("B130").Value = -14.53749430033371
$ws.Range("B131").Value = -14.72517126738194
$ws.Range("B132").Value = -9.684117434562555
$ws.Range("B133").Value = -2.675877847677288
$ws.Range("B134").Value = 3.854732564758434
$ws.Range("B135").Value = 7.057692386189451
$ws.Range("B136").Value = 10.81918261889252
$ws.Range("B137").Value = 12.7827260065003
$ws.Range("B138").Value = 10.472861510798
$ws.Range("B139").Value = 6.604385784822875
$ws.Range("B140").Value = 2.774439676894236
$ws.Range("B141").Value = 1.13146063314058
$ws.Range("B142").Value = 6.524097905960424
$ws.Range("B143").Value = 11.49679704207163
$ws.Range("B144").Value = 8.879514058481108
$ws.Range("B145").Value = 8.2355758601121
$ws.Range("B146").Value = 7.54384126261074
$ws.Range("B147").Value = 5.810477318753815
$ws.Range("B148").Value = 8.65272379870281
$ws.Range("B149").Value = 9.672573006707918
$ws.Range("B150").Value = 8.378752282192981
$ws.Range("B151").Value = 10.0805081826649
$ws.Range("B152").Value = 8.326592930107291
$ws.Range("B153").Value = 8.963830345108704
$ws.Range("B154").Value = 7.141668838797188
$ws.Range("B155").Value = 5.707693267442671
$ws.Range("B156").Value = 5.554947773991259
$ws.Range("B157").Value = -1.884082632649142
$ws.Range("B158").Value = -5.505710825134513
$ws.Range("B159").Value = -8.643826455346048
$ws.Range("B160").Value = -6.715623913875712
$ws.Range("B161").Value = 0.1940722961908836
$ws.Range("B162").Value = 5.730373278002543
$ws.Range("B163").Value = 12.37429431326815
$ws.Range("B164").Value = 10.35168698449283
$ws.Range("B165").Value = 8.299203222351625
$ws.Range("B166").Value = 7.990633317104523
$ws.Range("B167").Value = 5.55410427899734
$ws.Range("B168").Value = 6.262752151328743
$ws.Range("B169").Value = 4.900377161276204
$ws.Range("B170").Value = 1.613185587764487
$ws.Range("B171").Value = -1.470044815491289
$ws.Range("B172").Value = -1.196268114679577
$ws.Range("B173").Value = 0.485959315400164
$ws.Range("B174").Value = 0.7563940201671171
$ws.Range("B175").Value = 2.971778495402133
$ws.Range("B176").Value = 2.520277359953171
$ws.Range("B177").Value = 0.07307638773925706
$ws.Range("B178").Value = -0.9121082802822866
$ws.Range("B179").Value = -2.372371669152928
$ws.Range("B180").Value = -4.042059073304049
$ws.Range("B181").Value = -2.60821148365149
$ws.Range("B182").Value = 0.6089292777730293
$ws.Range("B183").Value = 2.607527860042458
$ws.Range("B184").Value = 4.178177039592734
$ws.Range("B185").Value = 2.491600889735968
$ws.Range("B186").Value = -0.4165027695021473
$ws.Range("B187").Value = -2.223860579068115
$ws.Range("B188").Value = -3.18347481753154
$ws.Range("B189").Value = -1.564967481364465
$ws.Range("B190").Value = 0.1244080080201826
$ws.Range("B191").Value = 3.000036687412866
$ws.Range("B192").Value = 3.891832179208521
$ws.Range("B193").Value = 4.159979339436637
$ws.Range("B194").Value = 3.244247247805476
$ws.Range("B195").Value = -3.530740872246796
$ws.Range("B196").Value = -3.51462917136951
$ws.Range("B197").Value = -6.547890298283832
$ws.Range("B198").Value = -5.791825167044163
$ws.Range("B199").Value = 1.020556135384498
$ws.Range("B200").Value = -1.65566852654232
$ws.Range("B201").Value = -1.385141389983435
$ws.Range("B202").Value = -7.962037462522275
$ws.Range("B203").Value = -13.84659963813887
$ws.Range("B204").Value = -10.2489542331945
$ws.Range("B205").Value = -4.195739289034717
$ws.Range("B206").Value = 6.407663956409415
$ws.Range("B207").Value = 11.95905640773656
$ws.Range("B208").Value = 11.47305607824434
$ws.Range("B209").Value = 8.416698832259083
$ws.Range("B210").Value = 6.627614743456678
$ws.Range("B211").Value = 8.763428670215244
$ws.Range("B212").Value = 6.212127479181206
$ws.Range("B213").Value = 2.530026725718892
$ws.Range("B214").Value = 0.2126104948851681
$ws.Range("B215").Value = -3.481583191854853
$ws.Range("B216").Value = -2.899131725008852
$ws.Range("B217").Value = -1.806096967880377
$ws.Range("B218").Value = -3.899549816607895
$ws.Range("B219").Value = -3.278656524395251
$ws.Range("B220").Value = -0.09459653849913652
$ws.Range("B221").Value = 2.516615113728693
$ws.Range("B222").Value = 5.711320481685123
$ws.Range("B223").Value = 6.825038461568234
$ws.Range("B224").Value = 3.188245782930421

